$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix text fields where commas were used as separators between names (should be periods)
$ws.Range("E91").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E186").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E199").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# Fix floating point formatting in the Importe column (H2:H261):
# remove thousands separators ("."), convert decimal comma (",") to decimal point (".")
# Force the range to remain text (it was stored as text before, e.g. "5.050,00")
$importeRange = $ws.Range("H2:H261")
$importeRange.NumberFormat = "@"

$ws.Range("H2").Value = "5050.00"
$ws.Range("H3").Value = "74280.00"
$ws.Range("H4").Value = "49200.00"
$ws.Range("H5").Value = "70940.00"
$ws.Range("H6").Value = "374220.00"
$ws.Range("H7").Value = "33600.00"
$ws.Range("H8").Value = "38000.00"
$ws.Range("H9").Value = "639572.00"
$ws.Range("H10").Value = "68000.00"
$ws.Range("H11").Value = "287334.00"
$ws.Range("H12").Value = "472.90"
$ws.Range("H13").Value = "2400.00"
$ws.Range("H14").Value = "275.76"
$ws.Range("H15").Value = "1280.00"
$ws.Range("H16").Value = "950.00"
$ws.Range("H17").Value = "773405.08"
$ws.Range("H18").Value = "456484.03"
$ws.Range("H19").Value = "92636.00"
$ws.Range("H20").Value = "182423.73"
$ws.Range("H21").Value = "10040.00"
$ws.Range("H22").Value = "32270.00"
$ws.Range("H23").Value = "25944.82"
$ws.Range("H24").Value = "16819.00"
$ws.Range("H25").Value = "23887.75"
$ws.Range("H26").Value = "10805.00"
$ws.Range("H27").Value = "25190.52"
$ws.Range("H28").Value = "25216.00"
$ws.Range("H29").Value = "32800.00"
$ws.Range("H30").Value = "3900.00"
$ws.Range("H31").Value = "3900.00"
$ws.Range("H32").Value = "1050.00"
$ws.Range("H33").Value = "1500.00"
$ws.Range("H34").Value = "3550.00"
$ws.Range("H35").Value = "3000.00"
$ws.Range("H36").Value = "11550.00"
$ws.Range("H37").Value = "20000.00"
$ws.Range("H38").Value = "25650.00"
$ws.Range("H39").Value = "79191.00"
$ws.Range("H40").Value = "420.00"
$ws.Range("H41").Value = "35500.00"
$ws.Range("H42").Value = "18218.53"
$ws.Range("H43").Value = "9588.00"
$ws.Range("H44").Value = "1118.59"
$ws.Range("H45").Value = "1189.50"
$ws.Range("H46").Value = "3928.07"
$ws.Range("H47").Value = "8399.99"
$ws.Range("H48").Value = "326390.31"
$ws.Range("H49").Value = "21540.19"
$ws.Range("H50").Value = "312.40"
$ws.Range("H51").Value = "1698.90"
$ws.Range("H52").Value = "2625.13"
$ws.Range("H53").Value = "570000.00"
$ws.Range("H54").Value = "8050.00"
$ws.Range("H55").Value = "14602.00"
$ws.Range("H56").Value = "4542.74"
$ws.Range("H57").Value = "8527.29"
$ws.Range("H58").Value = "2905.52"
$ws.Range("H59").Value = "126.00"
$ws.Range("H60").Value = "102.53"
$ws.Range("H61").Value = "69712.82"
$ws.Range("H62").Value = "200.00"
$ws.Range("H63").Value = "40224.18"
$ws.Range("H64").Value = "450.00"
$ws.Range("H65").Value = "1470.00"
$ws.Range("H66").Value = "2532.86"
$ws.Range("H67").Value = "852.00"
$ws.Range("H68").Value = "5760.00"
$ws.Range("H69").Value = "8730.00"
$ws.Range("H70").Value = "29764.00"
$ws.Range("H71").Value = "30960.00"
$ws.Range("H72").Value = "324800.00"
$ws.Range("H73").Value = "469.84"
$ws.Range("H74").Value = "3736.00"
$ws.Range("H75").Value = "835.00"
$ws.Range("H76").Value = "1220.00"
$ws.Range("H77").Value = "6230.00"
$ws.Range("H78").Value = "10693.00"
$ws.Range("H79").Value = "6400.00"
$ws.Range("H80").Value = "2820.00"
$ws.Range("H81").Value = "4160.00"
$ws.Range("H82").Value = "650.00"
$ws.Range("H83").Value = "605.00"
$ws.Range("H84").Value = "7100.40"
$ws.Range("H85").Value = "16200.00"
$ws.Range("H86").Value = "25230.00"
$ws.Range("H87").Value = "10.98"
$ws.Range("H88").Value = "1690.00"
$ws.Range("H89").Value = "1385.00"
$ws.Range("H90").Value = "44486.20"
$ws.Range("H91").Value = "1710.00"
$ws.Range("H92").Value = "6930.00"
$ws.Range("H93").Value = "9080.00"
$ws.Range("H94").Value = "2759.00"
$ws.Range("H95").Value = "526097.59"
$ws.Range("H96").Value = "1050.00"
$ws.Range("H97").Value = "150.00"
$ws.Range("H98").Value = "53261.49"
$ws.Range("H99").Value = "4000.00"
$ws.Range("H100").Value = "2500.00"
$ws.Range("H101").Value = "3450.00"
$ws.Range("H102").Value = "70.27"
$ws.Range("H103").Value = "79.72"
$ws.Range("H104").Value = "28238.58"
$ws.Range("H105").Value = "390.00"
$ws.Range("H106").Value = "185.00"
$ws.Range("H107").Value = "1749.00"
$ws.Range("H108").Value = "3633.00"
$ws.Range("H109").Value = "1905.10"
$ws.Range("H110").Value = "12046.60"
$ws.Range("H111").Value = "5542.08"
$ws.Range("H112").Value = "360.00"
$ws.Range("H113").Value = "1880.00"
$ws.Range("H114").Value = "14481.05"
$ws.Range("H115").Value = "831.00"
$ws.Range("H116").Value = "89617.99"
$ws.Range("H117").Value = "4753.39"
$ws.Range("H118").Value = "7570.50"
$ws.Range("H119").Value = "909.24"
$ws.Range("H120").Value = "440.00"
$ws.Range("H121").Value = "341.44"
$ws.Range("H122").Value = "278.70"
$ws.Range("H123").Value = "49200.54"
$ws.Range("H124").Value = "4218.01"
$ws.Range("H125").Value = "7131.50"
$ws.Range("H126").Value = "8441.76"
$ws.Range("H127").Value = "328.60"
$ws.Range("H128").Value = "1505.00"
$ws.Range("H129").Value = "12700.00"
$ws.Range("H130").Value = "2650.00"
$ws.Range("H131").Value = "613.80"
$ws.Range("H132").Value = "378.00"
$ws.Range("H133").Value = "7990.00"
$ws.Range("H134").Value = "600.00"
$ws.Range("H135").Value = "60434.00"
$ws.Range("H136").Value = "45386.00"
$ws.Range("H137").Value = "6033.00"
$ws.Range("H138").Value = "1845.00"
$ws.Range("H139").Value = "1738.50"
$ws.Range("H140").Value = "4820.00"
$ws.Range("H141").Value = "8782.00"
$ws.Range("H142").Value = "950.00"
$ws.Range("H143").Value = "18000.00"
$ws.Range("H144").Value = "3000.00"
$ws.Range("H145").Value = "4500.00"
$ws.Range("H146").Value = "5000.00"
$ws.Range("H147").Value = "131390.00"
$ws.Range("H148").Value = "1495.00"
$ws.Range("H149").Value = "1784.54"
$ws.Range("H150").Value = "129.00"
$ws.Range("H151").Value = "200.00"
$ws.Range("H152").Value = "2733.28"
$ws.Range("H153").Value = "7635.60"
$ws.Range("H154").Value = "2017.36"
$ws.Range("H155").Value = "5639.64"
$ws.Range("H156").Value = "15320.00"
$ws.Range("H157").Value = "8000.00"
$ws.Range("H158").Value = "3500.00"
$ws.Range("H159").Value = "10000.00"
$ws.Range("H160").Value = "26716.80"
$ws.Range("H161").Value = "2500.00"
$ws.Range("H162").Value = "3204.50"
$ws.Range("H163").Value = "2556.00"
$ws.Range("H164").Value = "3000.00"
$ws.Range("H165").Value = "3000.00"
$ws.Range("H166").Value = "2000.00"
$ws.Range("H167").Value = "1500.00"
$ws.Range("H168").Value = "17874.50"
$ws.Range("H169").Value = "7500.00"
$ws.Range("H170").Value = "4000.00"
$ws.Range("H171").Value = "2500.00"
$ws.Range("H172").Value = "4100.00"
$ws.Range("H173").Value = "11000.00"
$ws.Range("H174").Value = "15000.00"
$ws.Range("H175").Value = "6000.00"
$ws.Range("H176").Value = "2500.00"
$ws.Range("H177").Value = "3000.00"
$ws.Range("H178").Value = "27240.00"
$ws.Range("H179").Value = "3400.00"
$ws.Range("H180").Value = "900.00"
$ws.Range("H181").Value = "500.00"
$ws.Range("H182").Value = "1400.00"
$ws.Range("H183").Value = "10700.00"
$ws.Range("H184").Value = "6397.00"
$ws.Range("H185").Value = "964.02"
$ws.Range("H186").Value = "4622.00"
$ws.Range("H187").Value = "10638.00"
$ws.Range("H188").Value = "322.00"
$ws.Range("H189").Value = "7750.00"
$ws.Range("H190").Value = "12720.00"
$ws.Range("H191").Value = "1720.00"
$ws.Range("H192").Value = "1029.92"
$ws.Range("H193").Value = "22400.20"
$ws.Range("H194").Value = "1182.00"
$ws.Range("H195").Value = "1337.18"
$ws.Range("H196").Value = "1308.00"
$ws.Range("H197").Value = "46670.70"
$ws.Range("H198").Value = "28285.00"
$ws.Range("H199").Value = "4650.00"
$ws.Range("H200").Value = "7644.30"
$ws.Range("H201").Value = "14112.50"
$ws.Range("H202").Value = "5200.00"
$ws.Range("H203").Value = "4280.00"
$ws.Range("H204").Value = "18905.00"
$ws.Range("H205").Value = "2886.76"
$ws.Range("H206").Value = "4685.00"
$ws.Range("H207").Value = "5600.00"
$ws.Range("H208").Value = "25000.00"
$ws.Range("H209").Value = "25000.00"
$ws.Range("H210").Value = "8500.00"
$ws.Range("H211").Value = "25000.00"
$ws.Range("H212").Value = "25000.00"
$ws.Range("H213").Value = "50000.00"
$ws.Range("H214").Value = "50000.00"
$ws.Range("H215").Value = "25000.00"
$ws.Range("H216").Value = "6000.00"
$ws.Range("H217").Value = "4522.62"
$ws.Range("H218").Value = "216.86"
$ws.Range("H219").Value = "25900.00"
$ws.Range("H220").Value = "105000.00"
$ws.Range("H221").Value = "127500.00"
$ws.Range("H222").Value = "105000.00"
$ws.Range("H223").Value = "105000.00"
$ws.Range("H224").Value = "105000.00"
$ws.Range("H225").Value = "105000.00"
$ws.Range("H226").Value = "175000.00"
$ws.Range("H227").Value = "175000.00"
$ws.Range("H228").Value = "269000.00"
$ws.Range("H229").Value = "105000.00"
$ws.Range("H230").Value = "105000.00"
$ws.Range("H231").Value = "105000.00"
$ws.Range("H232").Value = "105000.00"
$ws.Range("H233").Value = "105000.00"
$ws.Range("H234").Value = "175000.00"
$ws.Range("H235").Value = "333000.00"
$ws.Range("H236").Value = "175000.00"
$ws.Range("H237").Value = "105000.00"
$ws.Range("H238").Value = "155000.00"
$ws.Range("H239").Value = "105000.00"
$ws.Range("H240").Value = "105000.00"
$ws.Range("H241").Value = "105000.00"
$ws.Range("H242").Value = "146481.60"
$ws.Range("H243").Value = "49968.24"
$ws.Range("H244").Value = "17143.30"
$ws.Range("H245").Value = "34880.00"
$ws.Range("H246").Value = "14000.00"
$ws.Range("H247").Value = "16899.00"
$ws.Range("H248").Value = "41967.28"
$ws.Range("H249").Value = "2500.00"
$ws.Range("H250").Value = "415.64"
$ws.Range("H251").Value = "7500.00"
$ws.Range("H252").Value = "7500.00"
$ws.Range("H253").Value = "20540.18"
$ws.Range("H254").Value = "6000.00"
$ws.Range("H255").Value = "1500.00"
$ws.Range("H256").Value = "7000.00"
$ws.Range("H257").Value = "9680.00"
$ws.Range("H258").Value = "44600.00"
$ws.Range("H259").Value = "651900.00"
$ws.Range("H260").Value = "22300.00"
$ws.Range("H261").Value = "295.00"

# Restore default (unstyled) cell style now that the text values are committed
$importeRange.Style = "Normal"

